$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 18.82499266666667
$ws.Range("H2").Value = 56.474978
$ws.Range("I2").Value = 0.06886869772378311
$ws.Range("J2").Value = 0.0688686977237831
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 8.131233999999999
$ws.Range("N2").Value = 24.393702
$ws.Range("O2").Value = 0.02090995573015822
$ws.Range("P2").Value = 0.02090995573015823
$ws.Range("Q2").Value = 153.0704204209507
$ws.Range("R2").Value = 1377.633783788556
$ws.Range("S2").Value = 0.001440041420597953
$ws.Range("T2").Value = 0.001440041420597953

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 18.82499266666667
$ws.Range("H3").Value = 56.474978
$ws.Range("I3").Value = 0.06886869772378311
$ws.Range("J3").Value = 0.0688686977237831
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 243.3763986666667
$ws.Range("N3").Value = 730.1291960000001
$ws.Range("O3").Value = 0.625857000534647
$ws.Range("P3").Value = 0.6258570005346471
$ws.Range("Q3").Value = 4581.558920139743
$ws.Range("R3").Value = 41234.03028125769
$ws.Range("S3").Value = 0.04310195658813417
$ws.Range("T3").Value = 0.04310195658813416

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 18.82499266666667
$ws.Range("H4").Value = 56.474978
$ws.Range("I4").Value = 0.06886869772378311
$ws.Range("J4").Value = 0.0688686977237831
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 103.9426383333333
$ws.Range("N4").Value = 311.827915
$ws.Range("O4").Value = 0.2672947262403034
$ws.Range("P4").Value = 0.2672947262403035
$ws.Range("Q4").Value = 1956.719404378986
$ws.Range("R4").Value = 17610.47463941087
$ws.Range("S4").Value = 0.01840823970460482
$ws.Range("T4").Value = 0.01840823970460482

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 18.82499266666667
$ws.Range("H5").Value = 56.474978
$ws.Range("I5").Value = 0.06886869772378311
$ws.Range("J5").Value = 0.0688686977237831
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 33.41874933333333
$ws.Range("N5").Value = 100.256248
$ws.Range("O5").Value = 0.08593831749489127
$ws.Range("P5").Value = 0.08593831749489128
$ws.Range("Q5").Value = 629.1077111291715
$ws.Range("R5").Value = 5661.969400162544
$ws.Range("S5").Value = 0.005918460010446169
$ws.Range("T5").Value = 0.005918460010446168

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 121.8208923333333
$ws.Range("H6").Value = 365.462677
$ws.Range("I6").Value = 0.4456653109566078
$ws.Range("J6").Value = 0.4456653109566078
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 8.131233999999999
$ws.Range("N6").Value = 24.393702
$ws.Range("O6").Value = 0.02090995573015822
$ws.Range("P6").Value = 0.02090995573015823
$ws.Range("Q6").Value = 990.5541816511394
$ws.Range("R6").Value = 8914.987634860254
$ws.Range("S6").Value = 0.009318841922569867
$ws.Range("T6").Value = 0.009318841922569869

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 121.8208923333333
$ws.Range("H7").Value = 365.462677
$ws.Range("I7").Value = 0.4456653109566078
$ws.Range("J7").Value = 0.4456653109566078
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 243.3763986666667
$ws.Range("N7").Value = 730.1291960000001
$ws.Range("O7").Value = 0.625857000534647
$ws.Range("P7").Value = 0.6258570005346471
$ws.Range("Q7").Value = 29648.33005844642
$ws.Range("R7").Value = 266834.9705260178
$ws.Range("S7").Value = 0.2789227547576433
$ws.Range("T7").Value = 0.2789227547576433

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 121.8208923333333
$ws.Range("H8").Value = 365.462677
$ws.Range("I8").Value = 0.4456653109566078
$ws.Range("J8").Value = 0.4456653109566078
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 103.9426383333333
$ws.Range("N8").Value = 311.827915
$ws.Range("O8").Value = 0.2672947262403034
$ws.Range("P8").Value = 0.2672947262403035
$ws.Range("Q8").Value = 12662.38495324761
$ws.Range("R8").Value = 113961.4645792285
$ws.Range("S8").Value = 0.1191239872869462
$ws.Range("T8").Value = 0.1191239872869462

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 121.8208923333333
$ws.Range("H9").Value = 365.462677
$ws.Range("I9").Value = 0.4456653109566078
$ws.Range("J9").Value = 0.4456653109566078
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 33.41874933333333
$ws.Range("N9").Value = 100.256248
$ws.Range("O9").Value = 0.08593831749489127
$ws.Range("P9").Value = 0.08593831749489128
$ws.Range("Q9").Value = 4071.101864450655
$ws.Range("R9").Value = 36639.9167800559
$ws.Range("S9").Value = 0.0382997269894484
$ws.Range("T9").Value = 0.03829972698944841

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 87.673585
$ws.Range("H10").Value = 263.020755
$ws.Range("I10").Value = 0.3207419907481189
$ws.Range("J10").Value = 0.3207419907481188
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 8.131233999999999
$ws.Range("N10").Value = 24.393702
$ws.Range("O10").Value = 0.02090995573015822
$ws.Range("P10").Value = 0.02090995573015823
$ws.Range("Q10").Value = 712.8944352538899
$ws.Range("R10").Value = 6416.04991728501
$ws.Range("S10").Value = 0.006706700827345984
$ws.Range("T10").Value = 0.006706700827345984

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 87.673585
$ws.Range("H11").Value = 263.020755
$ws.Range("I11").Value = 0.3207419907481189
$ws.Range("J11").Value = 0.3207419907481188
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 243.3763986666667
$ws.Range("N11").Value = 730.1291960000001
$ws.Range("O11").Value = 0.625857000534647
$ws.Range("P11").Value = 0.6258570005346471
$ws.Range("Q11").Value = 21337.68137549589
$ws.Range("R11").Value = 192039.132379463
$ws.Range("S11").Value = 0.2007386202751292
$ws.Range("T11").Value = 0.2007386202751292

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 87.673585
$ws.Range("H12").Value = 263.020755
$ws.Range("I12").Value = 0.3207419907481189
$ws.Range("J12").Value = 0.3207419907481188
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 103.9426383333333
$ws.Range("N12").Value = 311.827915
$ws.Range("O12").Value = 0.2672947262403034
$ws.Range("P12").Value = 0.2672947262403035
$ws.Range("Q12").Value = 9113.023737041758
$ws.Range("R12").Value = 82017.21363337584
$ws.Range("S12").Value = 0.08573264261078838
$ws.Range("T12").Value = 0.08573264261078838

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 87.673585
$ws.Range("H13").Value = 263.020755
$ws.Range("I13").Value = 0.3207419907481189
$ws.Range("J13").Value = 0.3207419907481188
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 33.41874933333333
$ws.Range("N13").Value = 100.256248
$ws.Range("O13").Value = 0.08593831749489127
$ws.Range("P13").Value = 0.08593831749489128
$ws.Range("Q13").Value = 2929.941560269693
$ws.Range("R13").Value = 26369.47404242724
$ws.Range("S13").Value = 0.02756402703485532
$ws.Range("T13").Value = 0.02756402703485532

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 45.02666966666666
$ws.Range("H14").Value = 135.080009
$ws.Range("I14").Value = 0.1647240005714903
$ws.Range("J14").Value = 0.1647240005714903
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 8.131233999999999
$ws.Range("N14").Value = 24.393702
$ws.Range("O14").Value = 0.02090995573015822
$ws.Range("P14").Value = 0.02090995573015823
$ws.Range("Q14").Value = 366.1223873003686
$ws.Range("R14").Value = 3295.101485703317
$ws.Range("S14").Value = 0.00344437155964442
$ws.Range("T14").Value = 0.00344437155964442

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 45.02666966666666
$ws.Range("H15").Value = 135.080009
$ws.Range("I15").Value = 0.1647240005714903
$ws.Range("J15").Value = 0.1647240005714903
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 243.3763986666667
$ws.Range("N15").Value = 730.1291960000001
$ws.Range("O15").Value = 0.625857000534647
$ws.Range("P15").Value = 0.6258570005346471
$ws.Range("Q15").Value = 10958.42870742697
$ws.Range("R15").Value = 98625.85836684277
$ws.Range("S15").Value = 0.1030936689137404
$ws.Range("T15").Value = 0.1030936689137404

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 45.02666966666666
$ws.Range("H16").Value = 135.080009
$ws.Range("I16").Value = 0.1647240005714903
$ws.Range("J16").Value = 0.1647240005714903
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 103.9426383333333
$ws.Range("N16").Value = 311.827915
$ws.Range("O16").Value = 0.2672947262403034
$ws.Range("P16").Value = 0.2672947262403035
$ws.Range("Q16").Value = 4680.190840516803
$ws.Range("R16").Value = 42121.71756465123
$ws.Range("S16").Value = 0.04402985663796408
$ws.Range("T16").Value = 0.04402985663796408

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 45.02666966666666
$ws.Range("H17").Value = 135.080009
$ws.Range("I17").Value = 0.1647240005714903
$ws.Range("J17").Value = 0.1647240005714903
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 33.41874933333333
$ws.Range("N17").Value = 100.256248
$ws.Range("O17").Value = 0.08593831749489127
$ws.Range("P17").Value = 0.08593831749489128
$ws.Range("Q17").Value = 1504.734986905137
$ws.Range("R17").Value = 13542.61488214623
$ws.Range("S17").Value = 0.01415610346014138
$ws.Range("T17").Value = 0.01415610346014138
